$d = $word.ActiveDocument

function ReplaceText($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# --- Hunk 1 (around JavaMethodService / EvaluationServices / AstEvaluator / AstSwitch) ---

ReplaceText `
    "`tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)`n`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)" `
    "`tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)`n`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)"

ReplaceText `
    "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:204)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:192)`n`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)" `
    "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:208)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)`n`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)"

ReplaceText `
    "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)" `
    "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)"

# --- Hunk 2 (GeneratedMethodAccessor) ---

ReplaceText `
    "`tat sun.reflect.GeneratedMethodAccessor74.invoke(Unknown Source)" `
    "`tat sun.reflect.GeneratedMethodAccessor73.invoke(Unknown Source)"

# --- Hunk 3 (tail of the stack trace, Maven/Tycho/Equinox block replaced by Eclipse JDT runner block) ---

$oldTail = "`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)`n" + `
    "`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)`n" + `
    "`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" + `
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" + `
    "`tat java.lang.reflect.Method.invoke(Method.java:498)`n" + `
    "`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)`n" + `
    "`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)`n" + `
    "`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)`n" + `
    "`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)`n" + `
    "`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" + `
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" + `
    "`tat java.lang.reflect.Method.invoke(Method.java:498)`n" + `
    "`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)`n" + `
    "`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)`n" + `
    "`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)`n" + `
    "`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)`n" + `
    "`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)`n" + `
    "`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" + `
    "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" + `
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" + `
    "`tat java.lang.reflect.Method.invoke(Method.java:498)`n" + `
    "`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)`n" + `
    "`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)`n" + `
    "`tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)`n" + `
    "`tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)"

$newTail = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n" + `
    "`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n" + `
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n" + `
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n" + `
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n" + `
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"

ReplaceText $oldTail $newTail

Write-Host "Done."
